# Add season record columns (Wins, Losses, Ties) to the SFG_2002 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - styled like other header cells (bold/border/center via style index 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the format of an existing header cell onto the new header cells so
# they pick up the same cell style (bold, bordered, centered) used by the
# rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-40: Wins=95, Losses=66, Ties=1 for every player row
for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 95
    $ws.Cells.Item($row, 31).Value = 66
    $ws.Cells.Item($row, 32).Value = 1
}
